$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows above row 12 (existing rows 12-23 shift down to 15-26)
$ws.Rows(12).Insert()
$ws.Rows(12).Insert()
$ws.Rows(12).Insert()

# New data for rows 12-14: Comercializadora del Agro de Limarí - Ciruela - Angeleno
# Boilerplate columns (A,B,C,E,F,G,H,I,J) are identical across all rows in this sheet.
$newRows = @(
    @{ Row = 12; Date = 44622; Variedad = "Angeleno"; Calidad = "Especial"; Volumen = 16;  Min = 235000; Max = 240000; Prom = 237500; Unidad = "$/bins (450 kilos)"; Origen = "Región Metropolitana"; PrecioKg = 528; KgUnidad = 450 },
    @{ Row = 13; Date = 44622; Variedad = "Angeleno"; Calidad = "Primera";  Volumen = 20;  Min = 195000; Max = 200000; Prom = 197500; Unidad = "$/bins (450 kilos)"; Origen = "Región Metropolitana"; PrecioKg = 439; KgUnidad = 450 },
    @{ Row = 14; Date = 44622; Variedad = "Angeleno"; Calidad = "Segunda";  Volumen = 20;  Min = 155000; Max = 160000; Prom = 157500; Unidad = "$/bins (450 kilos)"; Origen = "Región Metropolitana"; PrecioKg = 350; KgUnidad = 450 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = 2
    $ws.Cells.Item($row, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $r.Date
    $ws.Cells.Item($row, 5).Value = 4
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100103
    $ws.Cells.Item($row, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($row, 9).Value = 100103002
    $ws.Cells.Item($row, 10).Value = "Ciruela"
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
